$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 131-132, pushing existing rows 131.. down to 133..
$ws.Range("A131:T132").Insert()

# --- New row 131: Kurakata / Primera ---
$ws.Cells.Item(131, 1).Value  = 2
$ws.Cells.Item(131, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(131, 3).Value  = "Coquimbo"
$ws.Cells.Item(131, 4).Value  = 44931
$ws.Cells.Item(131, 5).Value  = 4
$ws.Cells.Item(131, 6).Value  = "Fruta"
$ws.Cells.Item(131, 7).Value  = 100103
$ws.Cells.Item(131, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(131, 9).Value  = 100103004
$ws.Cells.Item(131, 10).Value = "Durazno"
$ws.Cells.Item(131, 11).Value = "Kurakata"
$ws.Cells.Item(131, 12).Value = "Primera"
$ws.Cells.Item(131, 13).Value = 20
$ws.Cells.Item(131, 14).Value = 390000
$ws.Cells.Item(131, 15).Value = 400000
$ws.Cells.Item(131, 16).Value = 395000
$ws.Cells.Item(131, 17).Value = "$/bins (420 kilos)"
$ws.Cells.Item(131, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(131, 19).Value = 940
$ws.Cells.Item(131, 20).Value = 420

# --- New row 132: Kurakata / Segunda ---
$ws.Cells.Item(132, 1).Value  = 2
$ws.Cells.Item(132, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(132, 3).Value  = "Coquimbo"
$ws.Cells.Item(132, 4).Value  = 44931
$ws.Cells.Item(132, 5).Value  = 4
$ws.Cells.Item(132, 6).Value  = "Fruta"
$ws.Cells.Item(132, 7).Value  = 100103
$ws.Cells.Item(132, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(132, 9).Value  = 100103004
$ws.Cells.Item(132, 10).Value = "Durazno"
$ws.Cells.Item(132, 11).Value = "Kurakata"
$ws.Cells.Item(132, 12).Value = "Segunda"
$ws.Cells.Item(132, 13).Value = 20
$ws.Cells.Item(132, 14).Value = 330000
$ws.Cells.Item(132, 15).Value = 340000
$ws.Cells.Item(132, 16).Value = 335000
$ws.Cells.Item(132, 17).Value = "$/bins (420 kilos)"
$ws.Cells.Item(132, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(132, 19).Value = 798
$ws.Cells.Item(132, 20).Value = 420
